$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2025-07-23 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-07-24 Thursday", 2)

# Update the multiplication table cells by explicit (row, column) address so
# duplicate cell text (e.g. "41x79=3239" appears twice) is handled correctly.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "76×65=4940"
$t.Cell(1, 2).Range.Text  = "57×65=3705"
$t.Cell(1, 3).Range.Text  = "40×65=2600"
$t.Cell(1, 4).Range.Text  = "51×53=2703"
$t.Cell(1, 5).Range.Text  = "47×77=3619"

$t.Cell(5, 1).Range.Text  = "88×99=8712"
$t.Cell(5, 2).Range.Text  = "71×76=5396"
$t.Cell(5, 3).Range.Text  = "20×84=1680"
$t.Cell(5, 4).Range.Text  = "81×43=3483"
$t.Cell(5, 5).Range.Text  = "97×97=9409"

$t.Cell(10, 1).Range.Text = "91×33=3003"
$t.Cell(10, 2).Range.Text = "48×29=1392"
$t.Cell(10, 3).Range.Text = "22×34=748"
$t.Cell(10, 4).Range.Text = "62×50=3100"
$t.Cell(10, 5).Range.Text = "13×83=1079"

$t.Cell(15, 1).Range.Text = "97×44=4268"
$t.Cell(15, 2).Range.Text = "87×63=5481"
$t.Cell(15, 3).Range.Text = "58×12=696"
$t.Cell(15, 4).Range.Text = "46×35=1610"
$t.Cell(15, 5).Range.Text = "35×14=490"

$t.Cell(20, 1).Range.Text = "85×60=5100"
$t.Cell(20, 2).Range.Text = "94×85=7990"
$t.Cell(20, 3).Range.Text = "76×30=2280"
$t.Cell(20, 4).Range.Text = "22×26=572"
$t.Cell(20, 5).Range.Text = "23×20=460"
